# Auto-generated edit script: updates cryptos price/volume table
# to match the commit "Updated cryptos list on Sat Jun 22 14:13:08 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.282.97"
$ws.Range("E2").Value = "  +0.94%  "

$ws.Range("D3").Value = "3.504.09"
$ws.Range("E3").Value = "  +0.76%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.11%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.488"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.26%  "

$ws.Range("E9").Value = "  +2.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.389"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.79%  "

$ws.Range("D12").Value = "4.093.65"
$ws.Range("E12").Value = "  +0.91%  "

$ws.Range("E13").Value = "  +1.19%  "

$ws.Range("E14").Value = "  +2.72%  "

$ws.Range("D15").Value = "3.496.87"
$ws.Range("E15").Value = "  +0.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.16%  "

$ws.Range("D17").Value = "64.311.91"
$ws.Range("E17").Value = "  +0.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.570"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.30%  "

$ws.Range("D23").Value = "3.643.26"
$ws.Range("E23").Value = "  +0.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.39%  "

$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("E26").Value = "  +2.29%  "

$ws.Range("E27").Value = "  +2.43%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.21%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.92%  "

$ws.Range("E31").Value = "  +0.06%  "

$ws.Range("E32").Value = "  -5.64%  "

$ws.Range("D33").Value = "3.521.94"
$ws.Range("E33").Value = "  +1.18%  "

$ws.Range("E35").Value = "  +4.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.51%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "163.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0783"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.805"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.54%  "

$ws.Range("E47").Value = "  +1.86%  "

$ws.Range("E48").Value = "  -2.57%  "

$ws.Range("D49").Value = "2.478.19"
$ws.Range("E49").Value = "  +2.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.68%  "
